# Levelup curve limited to max 60 + show collection-unlocked ratio row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cap the level-up defeat threshold formula in C9 at 60.
$ws.Range("C9").Formula = "=MIN(60,FLOOR(1+(B9*B9*0.25),1))"

# Scroll/select so the "collection unlocked ratio" row (9) is visible,
# matching the new active cell/selection (was C11, now C9).
$ws.Range("C9").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
